# Archetypes and Aspects - wire up the newly-added archetypes/conditions.
# Commit message: "I believe this gets ticks working for Conditions."
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (Unarmed) ---
# "Battle Savant" is a new Pugilist sub-archetype; "Monk" (formerly in Q7)
# moves down into the new Nucleus/Monk row (row 13), so Q7 is cleared here.
$ws.Range("B7").Value = "Battle Savant"
$ws.Range("Q7").ClearContents()

# --- Row 11 (Sorcery) ---
# "Necromancer" shifts from H11 to K11, "Powderwitch" is inserted as a new
# Duskblade sub-archetype at C11 (old F11 "Necromancer" slot is cleared).
$ws.Range("C11").Value = "Powderwitch"
$ws.Range("F11").ClearContents()
$ws.Range("K11").Value = "Necromancer"
$ws.Range("H11").ClearContents()

# New blank "tick" cell for the Nimble condition on the Warlock row, styled
# distinctly (non-default format) but left without a value - mirrors the
# existing styled-but-empty convention used elsewhere (e.g. L1).
$ws.Range("N11").Style = "Normal"

# --- Row 12 (Sorcery) ---
# "Ardent" is renamed to "Dawnblade" (same cell/position, new name only).
$ws.Range("B12").Value = "Dawnblade"

# --- Row 13 (Mentalist) ---
# "Battle Savant" (old N13) and "Erased Spy" (old P13) both shift right to
# make room, and the old "Atom Lord" slot is renamed to "Nucleus".
$ws.Range("N13").ClearContents()
$ws.Range("P13").Value = "Erased Spy"
$ws.Range("Q13").Value = "Monk"
$ws.Range("V13").Value = "Nucleus"

# Restore the author's final cursor position on the new tick cell.
$ws.Range("N11").Select()
